$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely by the naive component forecaster bug fix
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recalculated forecast values (tiny floating point corrections)
$ws.Range("E3").Value = 0.4748521911469572
$ws.Range("E4").Value = 0.6970543652217387
$ws.Range("C6").Value = -0.01587181126745385
$ws.Range("C7").Value = 0.5978820435290855
$ws.Range("C8").Value = -0.02256889165886955
$ws.Range("E8").Value = -0.02753509623224515
$ws.Range("C9").Value = -0.50613598754502
$ws.Range("C10").Value = 0.09611428386595566
$ws.Range("C11").Value = -0.2706540469742613
$ws.Range("E12").Value = -0.2251688766574889
$ws.Range("E13").Value = 0.01247916696662799
$ws.Range("C14").Value = -0.001350220946472191
$ws.Range("E15").Value = 0.174086048246691
$ws.Range("E17").Value = -0.03768624985648339
$ws.Range("C18").Value = -0.5761528471665334
$ws.Range("C19").Value = -0.03047919532178645
$ws.Range("E19").Value = -0.1249617237519041
$ws.Range("C21").Value = 0.1932702877606163
$ws.Range("E21").Value = -0.2500935825088479
$ws.Range("E23").Value = 0.2001500500062203
$ws.Range("C27").Value = 0.2799548089016834
$ws.Range("E27").Value = 0.3604862916655849
$ws.Range("C29").Value = -0.2916219766884165
$ws.Range("C30").Value = -0.4278219446121501
$ws.Range("C32").Value = -0.2932081122163255
$ws.Range("E32").Value = -0.112644651861793
$ws.Range("C33").Value = -0.3404043877497931
$ws.Range("E33").Value = 0.01966269405897503
$ws.Range("C34").Value = -1.026566979837429
$ws.Range("E35").Value = -0.5986513494937395
$ws.Range("C36").Value = -0.02187747290984809
$ws.Range("E36").Value = -0.1415666278731686
$ws.Range("E40").Value = 0.9379151023484189
$ws.Range("C42").Value = 0.4636049209196802
$ws.Range("E42").Value = 0.2986939435938973
$ws.Range("E43").Value = 0.4006004000999486
$ws.Range("C46").Value = 0.6216390921348403
$ws.Range("E46").Value = -0.0776179936130994
$ws.Range("C47").Value = -0.3371636084037011
$ws.Range("C50").Value = -0.6768900623516871
$ws.Range("E52").Value = -0.338776212162295
$ws.Range("C53").Value = 1.132847588656238
$ws.Range("E53").Value = 0.2467479214459667
